# "Fixes to decks 1 and 5"
#
# Slide 3 ("Redux") and Slide 5 ("Flux") had their section titles swapped
# by mistake - the deck should read Flux, Flux, Redux, Redux (slides 2-5),
# grouping the Flux slides together followed by the Redux slides.
#
# Swap the title text on slide 3 and slide 5 so the titles land on the
# right slide.

$p = $ppt.ActivePresentation

$slide3 = $p.Slides.Item(3)
$slide5 = $p.Slides.Item(5)

$title3 = $slide3.Shapes.Item(1).TextFrame.TextRange
$title5 = $slide5.Shapes.Item(1).TextFrame.TextRange

$text3 = $title3.Text
$text5 = $title5.Text

$title3.Text = $text5
$title5.Text = $text3
